$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (row 1) ----
$ws.Range("A1").Value = "Author"
$ws.Range("B1").Value = "Type of data"
$ws.Range("C1").Value = "Processed?"
$ws.Range("D1").Value = "Normalisation"
$ws.Range("E1").Value = "Accessed from"

# ---- Data rows (2-13) ----
# Row 2
$ws.Range("A2").Value = "Clement-Ziza et al. (2014)"
$ws.Range("B2").Value = "RNA-seq of segregants grown on EMM"
$ws.Range("C2").Value = "Raw count matrix"
$ws.Range("D2").Value = "DESeq"
$ws.Range("E2").Value = "Supplementary dataset S4"

# Row 3
$ws.Range("A3").Value = "Clement-Ziza et al. (2014)"
$ws.Range("B3").Value = "Whole-genome sequencing of parental strains (JB50, JB759)"
$ws.Range("C3").Value = "Raw FASTQ"
$ws.Range("D3").Value = "NA"
$ws.Range("E3").Value = "ENA: ERX007392 and ERX007395"

# Row 4
$ws.Range("A4").Value = "Rubio et al. (2021)"
$ws.Range("B4").Value = "RNA-seq data in environmental stresses"
$ws.Range("C4").Value = "Raw count matrix"
$ws.Range("D4").Value = "DESeq"
$ws.Range("E4").Value = "Supplementary files"

# Row 5
$ws.Range("A5").Value = "Jeffares et al. (2015)"
$ws.Range("B5").Value = "Short variants in wild isolates"
$ws.Range("C5").Value = "Processed .vcf files"
$ws.Range("D5").Value = "NA"
$ws.Range("E5").Value = "Figshare associate with publication"

# Row 6
$ws.Range("A6").Value = "Atkinson et al. (2018)"
$ws.Range("B6").Value = "RNA-seq in various conditions (including EMM and YES)"
$ws.Range("C6").Value = "Processed gene-conditions fold-change matrix"
$ws.Range("D6").Value = "NA"
$ws.Range("E6").Value = "Supplementary table 3."

# Row 7
$ws.Range("A7").Value = "Kwon et al. (2012)"
$cell = $ws.Range("B7")
$cell.Value = "ChIP-seq data on mbx2-HA strain"
$cell.Characters(18, 4).Font.Italic = $true
$ws.Range("C7").Value = "Processed list of hits"
$ws.Range("D7").Value = "NA"
$ws.Range("E7").Value = "Supplementary files"

# Row 8
$ws.Range("A8").Value = "Kwon et al. (2012)"
$cell = $ws.Range("B8")
$cell.Value = "Microarray data on mbx2OE strain"
$cell.Characters(20, 4).Font.Italic = $true
$ws.Range("C8").Value = "Processed list of fold-change values"
$ws.Range("D8").Value = "NA"
$ws.Range("E8").Value = "Supplementary files"

# Row 9
$ws.Range("A9").Value = "Linder et al. (2008)"
$ws.Range("B9").Value = "Microarray data on Mediator deletion strains grown on EMM"
$ws.Range("C9").Value = "Processed matrix"
$ws.Range("D9").Value = "NA"
$ws.Range("E9").Value = "Supplementary files"

# Row 10
$ws.Range("A10").Value = "Szilagyi et al. (2012)"
$cell = $ws.Range("B10")
$cell.Value = "Microarray data on fkh2 deletion strain"
$cell.Characters(19, 5).Font.Italic = $true
$ws.Range("C10").Value = "Raw .CEL files"
$ws.Range("D10").Value = "MAS5, using Affy"
$ws.Range("E10").Value = "GEO: GSE31642"

# Row 11
$ws.Range("A11").Value = "Garg et al. (2015)"
$cell = $ws.Range("B11")
$cell.Value = "Microarray data on fkh2 deletion strain across timepoints"
$cell.Characters(20, 4).Font.Italic = $true
$ws.Range("C11").Value = "Processed fold-change matrix"
$ws.Range("D11").Value = "NA"
$ws.Range("E11").Value = "GEO: GSE60718"

# Row 12
$ws.Range("A12").Value = "Garg et al. (2015)"
$cell = $ws.Range("B12")
$cell.Value = "ChIP-seq data on fkh2-TAP strain"
$cell.Characters(18, 4).Font.Italic = $true
$ws.Range("C12").Value = "Processed list of targets"
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "Supplementary files"

# Row 13
$ws.Range("A13").Value = "Garg et al. (2023)"
$ws.Range("B13").Value = "RNA-seq of lab strain grown in Phosphate-minus media"
$ws.Range("C13").Value = "Processed DEseq comparison between timepoints"
$ws.Range("D13").Value = "DESeq"
$ws.Range("E13").Value = "GEO: GSE217953"

# ---- Formatting ----
# Data rows 2-13: full thin-box border, no fill (same visual style used for rows 2-11 previously).
$dataRange = $ws.Range("A2:E13")
$dataRange.Borders.LineStyle = 1

# Header row 1: thin box border, gray fill, but WITHOUT a bottom border (new in this revision).
$headerRange = $ws.Range("A1:E1")
$headerRange.Borders.LineStyle = 1
$headerRange.Interior.Color = 13553360
$headerRange.Borders(9).LineStyle = -4142

# ---- Selection ----
$ws.Range("B16").Select()
